$wb = $excel.ActiveWorkbook

# --- "definitions" sheet: update end-user disposal & end-user privacy rows
#     from "beneficial" to "non-beneficial" (category_binary formula in G
#     recalculates automatically: =IF(F="beneficial",1,0))
$wsDef = $wb.Worksheets.Item("definitions")
$wsDef.Range("F24").Value = "non-beneficial"
$wsDef.Range("F26").Value = "non-beneficial"

# --- "indicator_type" sheet: update technology scores (SI version 3)
$wsInd = $wb.Worksheets.Item("indicator_type")
$wsInd.Range("W3").Value = 0
$wsInd.Range("Y3").Value = 0

# Move the selection from I33 to F10
$wsInd.Activate()
$wsInd.Range("F10").Select()
